$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> updated numeric value, per the diff (imputed values
# recalculated for the RandomForest result, "Update Name of Algo").
$cellUpdates = @{
    "B2" = 8.420300000000003
    "A3" = -21.8546
    "C3" = -11.34810000000001
    "D6" = -7.884599999999995
    "E8" = 16.42210000000001
    "C12" = -11.82569999999999
    "A14" = -21.6437
    "D19" = -8.606799999999993
    "A21" = -20.36059999999999
    "A23" = -19.86239999999998
    "E23" = 16.31459999999999
    "C24" = -13.00419999999999
    "D24" = -7.978199999999995
    "A25" = -21.8123
    "B25" = 5.822799999999998
    "C25" = -12.99229999999999
    "A26" = -21.21339999999996
    "E26" = 16.397
    "B27" = 5.964600000000001
    "A29" = -20.67359999999998
    "D30" = -7.567499999999998
    "B31" = 5.074400000000001
    "D31" = -8.170399999999995
    "D33" = -7.696899999999996
    "E37" = 16.60400000000002
    "B39" = 9.952100000000003
    "D42" = -8.973699999999996
    "B48" = 5.428400000000003
    "E48" = 17.3589
    "C50" = -13.58389999999999
    "B51" = 5.531100000000001
    "B52" = 4.746600000000003
    "A53" = -21.5059
    "C53" = -10.4108
    "B55" = 5.690799999999997
    "D55" = -8.1852
    "B56" = 5.0914
    "A57" = -22.02270000000001
    "B57" = 5.364699999999996
    "C57" = -13.00879999999999
    "D58" = -8.102499999999997
    "A59" = -22.20819999999999
    "C61" = -12.62739999999999
    "E62" = 16.5658
    "C63" = -11.3061
    "D65" = -8.068599999999996
    "E66" = 17.00940000000002
    "A69" = -21.60799999999998
    "C70" = -11.9446
    "D70" = -8.160999999999998
    "B73" = 8.309499999999993
    "D75" = -8.031700000000003
    "A79" = -20.4513
    "A83" = -22.0772
    "D83" = -8.205999999999991
    "C86" = -13.39269999999999
    "D86" = -7.893399999999988
    "B89" = 4.875399999999997
    "E89" = 17.81870000000001
    "B90" = 5.4752
    "A91" = -21.26350000000001
    "B92" = 5.142899999999994
    "A93" = -20.81669999999998
    "E94" = 18.72990000000002
    "D96" = -8.004499999999993
    "D97" = -8.102
    "C98" = -12.02059999999999
    "C100" = -13.16549999999999
    "C102" = -13.07930000000001
}

foreach ($addr in $cellUpdates.Keys) {
    $ws.Range($addr).Value = $cellUpdates[$addr]
}
